# Update mods data [2025-12-21 15:09:08]
# Append a new row (42) to the ModCounts sheet with the latest mod count.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ModCounts")

$lastRow = 41
$newRow = 42

# Stamp the new row's formatting from the last existing row so it matches
# the sheet's established style (center-aligned, same font/fill/border).
$ws.Range("A" + $lastRow + ":C" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)

# Column A holds a date formatted as plain text (e.g. "2025/12/20"), not a
# real date value. Force the new cell to Text first so the date-like string
# isn't auto-converted into a date serial number, then re-stamp the format
# from the row above so the cell keeps the normal General/center style.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025/12/21"
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1097

$excel.CutCopyMode = 0
